$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, set its value then copy the formatting (bold header
# style) from the adjacent existing header cell G1 so it matches the rest
# of the header row.
$ws.Range("H1").Value = "ubicacion_descarga"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cells H2 and H3 (plain, unstyled like the rest of the data rows)
$ws.Range("H2").Value = "./descargas/RCEL/20987654321"
$ws.Range("H3").Value = "./descargas/RCEL/20999999999"
